$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShapeImage {
    param(
        $Shape,
        [string]$OldName,
        [string]$NewName
    )

    $r = $Shape.Range
    $xml = $r.WordOpenXML
    $oldAttr = 'name="' + $OldName + '"'
    $newAttr = 'name="' + $NewName + '"'
    $updated = $xml.Replace($oldAttr, $newAttr)
    $r.WordOpenXML = $updated
}

# Footer (default / primary) -> physical footer2.xml, docPr/cNvPr id="2", image1.png -> image2.png
$ftr1 = $sec.Footers.Item(1)
Rename-InlineShapeImage -Shape $ftr1.Range.InlineShapes.Item(1) -OldName "image1.png" -NewName "image2.png"

# Footer (first page) -> physical footer1.xml, docPr/cNvPr id="3", image1.png -> image2.png
$ftr2 = $sec.Footers.Item(2)
Rename-InlineShapeImage -Shape $ftr2.Range.InlineShapes.Item(1) -OldName "image1.png" -NewName "image2.png"

# Header (first page) -> physical header1.xml, docPr/cNvPr id="1", image2.jpg -> image1.jpg
$hdr2 = $sec.Headers.Item(2)
Rename-InlineShapeImage -Shape $hdr2.Range.InlineShapes.Item(1) -OldName "image2.jpg" -NewName "image1.jpg"

Write-Output "Done renaming inline shape image names"
